# Adapt tests to control version
# Adds a "version" column to the settings sheet, giving it the value 1,
# so the form carries a controllable version number.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("settings")

# New header cell for the version column. C1 already carries the sheet's
# style index 5 (it was a pre-existing blank styled cell), so just fill in
# its text value.
$ws.Range("C1").Value = "version"

# Numeric version value under the new header.
$ws.Range("C2").Value = 1

# Move the active selection, matching the post-edit worksheet state.
$ws.Range("C3").Select()
